$d = $word.ActiveDocument

# Helper: replace a paragraph's text while preserving the paragraph's own
# <w:pPr> and the leading empty <w:r/> run that this document consistently
# uses ahead of its "real" text run. We do this by inserting a brand-new
# paragraph (which inherits formatting from its neighbour), filling it via
# InsertXML (so literal apostrophes survive instead of being mangled by
# AutoCorrect's smart-quote substitution), and then deleting the original
# paragraph.
function Replace-ParaKeepLeadEmpty($doc, $oldText, $pPrXml, $runXml) {
    $idx = -1
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq $oldText) {
            $idx = $i
            break
        }
    }
    if ($idx -eq -1) {
        throw "Paragraph not found: $oldText"
    }

    $oldPara = $doc.Paragraphs.Item($idx)
    $oldRng = $oldPara.Range
    [void]$oldRng.Collapse(0)
    [void]$oldRng.InsertParagraphAfter()

    $newPara = $doc.Paragraphs.Item($idx + 1)
    $xmlFrag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$pPrXml<w:r/>$runXml</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    [void]$newPara.Range.InsertXML($xmlFrag)

    [void]$doc.Paragraphs.Item($idx).Range.Delete()
}

# --- 1. Title: shorten heading text (direct Range.Text assignment avoids
#          smart-quote autocorrect mangling the apostrophe) ---
$p1 = $d.Paragraphs.Item(1)
$rng1 = $p1.Range
$rng1.MoveEnd(1, -1)
$rng1.Text = "Play Lucky Lady's Charm for Free"

# --- 2. Drop the whole "Meta description: ..." paragraph that followed it ---
$d.Paragraphs.Item(2).Range.Delete()

# --- 3 & 4. "What we like" / "What we don't like" bullet list rewrites ---
$bulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

Replace-ParaKeepLeadEmpty $d "Autoplay feature for easier gameplay" $bulletPPr '<w:r><w:t>Flexible betting options</w:t></w:r>'
Replace-ParaKeepLeadEmpty $d "Intuitive interface and flexible betting options" $bulletPPr '<w:r><w:t>Intuitive interface</w:t></w:r>'
Replace-ParaKeepLeadEmpty $d "Free spins bonus round with tripled winnings" $bulletPPr '<w:r><w:t>Free spins bonus round</w:t></w:r>'
Replace-ParaKeepLeadEmpty $d "Playable on desktop, laptop, and mobile devices" $bulletPPr '<w:r><w:t>Available to play on desktop and mobile devices</w:t></w:r>'
Replace-ParaKeepLeadEmpty $d "Simple graphics and animations may not appeal to everyone" $bulletPPr '<w:r><w:t>Lack of flashy animations or sound effects</w:t></w:r>'
Replace-ParaKeepLeadEmpty $d "No progressive jackpot feature" $bulletPPr '<w:r><w:t>Limited number of pay lines</w:t></w:r>'

# --- 5. Insert a new bold paragraph "Play Lucky Lady's Charm for Free" right
#          after the "Limited number of pay lines" bullet, before the closing
#          (italic) paragraph. ---
$paraCount = $d.Paragraphs.Count
$bulletPara = $d.Paragraphs.Item($paraCount - 1)
$bulletRng = $bulletPara.Range
$bulletRng.Collapse(0)
$bulletRng.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($paraCount)
$xmlFrag2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Lucky Lady's Charm for Free</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$newPara.Range.InsertXML($xmlFrag2)

# --- 6. Replace the closing italic image-prompt paragraph with the new
#          closing sentence (direct Range.Text keeps the straight apostrophe
#          and preserves the existing italic run formatting) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRng = $lastPara.Range
$lastRng.MoveEnd(1, -1)
$lastRng.Text = "Read our review of Lucky Lady's Charm and play this nostalgic slot game for free."
